$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new place record as row 9
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Apartment Sayat-Nova 18"
$ws.Range("C9").Value = "Apartment"
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 40.181066166026902
$ws.Range("F9").Value = 44.521552090821501
$ws.Range("G9").Value = "Kentron"

# Update the active selection like the author's last interaction
$ws.Range("C19").Select()
